$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 02:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 738792
$ws.Range("C4").Value = 29057
$ws.Range("D4").Value = 68269
$ws.Range("E4").Value = 631509
$ws.Range("G4").Value = 1867
$ws.Range("H4").Value = 39014

# Row 8 - Alemania
$ws.Range("B8").Value = 143724
$ws.Range("C8").Value = 2327
$ws.Range("E8").Value = 53786
$ws.Range("G8").Value = 186
$ws.Range("H8").Value = 4538

# Row 38 - Australia
$ws.Range("B38").Value = 6586
$ws.Range("C38").Value = 53
$ws.Range("E38").Value = 2349

# Row 81 - Camerun
$ws.Range("D81").Value = 305
$ws.Range("E81").Value = 670
$ws.Range("F81").Value = 33
$ws.Range("G81").Value = 20
$ws.Range("H81").Value = 42

# Row 127 - Martinica
$ws.Range("B127").Value = 163
$ws.Range("C127").Value = 5
$ws.Range("E127").Value = 78
$ws.Range("F127").Value = 11
$ws.Range("G127").Value = 4
$ws.Range("H127").Value = 12

$wb.Save()
